$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B/C text updates (plain text, safe to set directly) ---
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("B14").Value = 'Solana'
$ws.Range("C14").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("B15").Value = 'BinanceUSD'
$ws.Range("C15").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("B38").Value = 'Algorand'
$ws.Range("C38").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("B39").Value = 'ARBITRUM'
$ws.Range("C39").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'

# --- Column D/E numeric-looking text updates staged via formula + PasteSpecial to avoid type coercion ---
$ws.Range("Z2").Formula = "=""27.990.33"""
$ws.Range("AA2").Formula = "=""  -0.26%  """
$ws.Range("Z3").Formula = "=""1.898.52"""
$ws.Range("AA3").Formula = "=""  +1.89%  """
$ws.Range("Z4").Formula = "=""1.003"""
$ws.Range("AA4").Formula = "=""  +0.13%  """
$ws.Range("Z5").Formula = "=""312.72"""
$ws.Range("AA5").Formula = "=""  +0.37%  """
$ws.Range("Z6").Formula = "=""1.003"""
$ws.Range("AA6").Formula = "=""  +0.17%  """
$ws.Range("Z7").Formula = "=""0.5018"""
$ws.Range("AA7").Formula = "=""  +0.78%  """
$ws.Range("Z8").Formula = "=""0.3921"""
$ws.Range("AA8").Formula = "=""  +0.47%  """
$ws.Range("Z9").Formula = "=""0.09499"""
$ws.Range("AA9").Formula = "=""  -2.41%  """
$ws.Range("Z10").Formula = "=""1.131"""
$ws.Range("AA10").Formula = "=""  -0.71%  """
$ws.Range("Z11").Formula = "=""41.96"""
$ws.Range("AA11").Formula = "=""  +2.81%  """
$ws.Range("Z12").Formula = "=""6.361"""
$ws.Range("AA12").Formula = "=""  -1.55%  """
$ws.Range("Z13").Formula = "=""1.906.75"""
$ws.Range("AA13").Formula = "=""  +2.31%  """
$ws.Range("Z14").Formula = "=""20.73"""
$ws.Range("AA14").Formula = "=""  -0.47%  """
$ws.Range("Z15").Formula = "=""1.003"""
$ws.Range("AA15").Formula = "=""  +0.15%  """
$ws.Range("Z16").Formula = "=""7.299"""
$ws.Range("AA16").Formula = "=""  -0.89%  """
$ws.Range("Z17").Formula = "=""0.00001115"""
$ws.Range("AA17").Formula = "=""  -0.74%  """
$ws.Range("Z18").Formula = "=""92.03"""
$ws.Range("AA18").Formula = "=""  -0.90%  """
$ws.Range("Z19").Formula = "=""0.06603"""
$ws.Range("AA19").Formula = "=""  +0.23%  """
$ws.Range("Z20").Formula = "=""17.79"""
$ws.Range("AA20").Formula = "=""  +2.05%  """
$ws.Range("Z21").Formula = "=""1.001"""
$ws.Range("AA21").Formula = "=""  -0.03%  """
$ws.Range("Z22").Formula = "=""6.187"""
$ws.Range("AA22").Formula = "=""  +1.24%  """
$ws.Range("Z23").Formula = "=""28.052.92"""
$ws.Range("AA23").Formula = "=""  -0.25%  """
$ws.Range("Z24").Formula = "=""11.26"""
$ws.Range("AA24").Formula = "=""  -0.43%  """
$ws.Range("Z25").Formula = "=""2.301"""
$ws.Range("AA25").Formula = "=""  +0.54%  """
$ws.Range("Z26").Formula = "=""2.649"""
$ws.Range("AA26").Formula = "=""  +4.29%  """
$ws.Range("Z27").Formula = "=""2.124.32"""
$ws.Range("AA27").Formula = "=""  +2.22%  """
$ws.Range("Z28").Formula = "=""20.76"""
$ws.Range("AA28").Formula = "=""  -1.41%  """
$ws.Range("Z29").Formula = "=""156.96"""
$ws.Range("AA29").Formula = "=""  +0.06%  """
$ws.Range("Z30").Formula = "=""126.82"""
$ws.Range("AA30").Formula = "=""  -0.41%  """
$ws.Range("Z31").Formula = "=""1.078"""
$ws.Range("AA31").Formula = "=""  +2.17%  """
$ws.Range("Z32").Formula = "=""0.1063"""
$ws.Range("AA32").Formula = "=""  +0.83%  """
$ws.Range("Z33").Formula = "=""5.601"""
$ws.Range("AA33").Formula = "=""  -0.22%  """
$ws.Range("Z34").Formula = "=""3.615"""
$ws.Range("AA34").Formula = "=""  -0.74%  """
$ws.Range("Z35").Formula = "=""9.597"""
$ws.Range("AA35").Formula = "=""  +1.84%  """
$ws.Range("Z36").Formula = "=""0.06580"""
$ws.Range("AA36").Formula = "=""  -2.25%  """
$ws.Range("Z37").Formula = "=""0.02424"""
$ws.Range("AA37").Formula = "=""  +1.36%  """
$ws.Range("Z38").Formula = "=""0.2172"""
$ws.Range("Z39").Formula = "=""1.224"""
$ws.Range("AA39").Formula = "=""  -0.66%  """
$ws.Range("Z40").Formula = "=""1.261"""
$ws.Range("AA40").Formula = "=""  +7.31%  """
$ws.Range("Z41").Formula = "=""4.975"""
$ws.Range("AA41").Formula = "=""  -0.32%  """
$ws.Range("Z42").Formula = "=""0.6316"""
$ws.Range("AA42").Formula = "=""  +0.79%  """
$ws.Range("Z43").Formula = "=""11.31"""
$ws.Range("AA43").Formula = "=""  -1.16%  """
$ws.Range("Z44").Formula = "=""1.003"""
$ws.Range("AA44").Formula = "=""  +0.13%  """
$ws.Range("Z45").Formula = "=""13.29"""
$ws.Range("AA45").Formula = "=""  -0.88%  """
$ws.Range("Z46").Formula = "=""0.5965"""
$ws.Range("AA46").Formula = "=""  -0.79%  """
$ws.Range("Z47").Formula = "=""3.719"""
$ws.Range("AA47").Formula = "=""  +1.89%  """
$ws.Range("Z48").Formula = "=""1.279"""
$ws.Range("AA48").Formula = "=""  +0.74%  """
$ws.Range("Z49").Formula = "=""2.012"""
$ws.Range("AA49").Formula = "=""  +1.91%  """
$ws.Range("Z50").Formula = "=""123.38"""
$ws.Range("AA50").Formula = "=""  -0.46%  """
$ws.Range("Z51").Formula = "=""1.176"""
$ws.Range("AA51").Formula = "=""  -1.40%  """

$ws.Range("Z2:AA51").Copy()
$ws.Range("D2").PasteSpecial(-4163, -4142, $true, $false)
$ws.Range("Z2:AA51").Clear()
$excel.CutCopyMode = $false
Write-Host "Done"